$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2735.125
$ws.Range("I86").Value = 1596.6
$ws.Range("J86").Value = 4632.6665
$ws.Range("K86").Value = 1596.6
$ws.Range("L86").Value = 4632.6665
$ws.Range("M86").Value = -473.5999999999999
$ws.Range("N86").Value = -6878.6665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 2735.125
$ws.Range("I89").Value = 1596.6
$ws.Range("J89").Value = 4632.6665
$ws.Range("K89").Value = 7983
$ws.Range("L89").Value = 23163.3325
$ws.Range("M89").Value = -2367
$ws.Range("N89").Value = -34395.3325

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 4202352
$ws.Range("I96").Value = 5494737.5
$ws.Range("J96").Value = 2099.75
$ws.Range("K96").Value = 16484212.5
$ws.Range("L96").Value = 6299.25
$ws.Range("M96").Value = -16482839.5
$ws.Range("N96").Value = -9045.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 716.7273
$ws.Range("I101").Value = 190.75
$ws.Range("J101").Value = 1017.2857
$ws.Range("K101").Value = 572.25
$ws.Range("L101").Value = 3051.8571
$ws.Range("M101").Value = 1049.75
$ws.Range("N101").Value = -6295.8571

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3339984.8
$ws.Range("I132").Value = 3590393.5
$ws.Range("J132").Value = 1199.6666
$ws.Range("K132").Value = 10771180.5
$ws.Range("L132").Value = 3598.9998
$ws.Range("M132").Value = -10768650.5
$ws.Range("N132").Value = -8658.9998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 22727.465
$ws.Range("I137").Value = 16570.63
$ws.Range("J137").Value = 35725.223
$ws.Range("K137").Value = 49711.89
$ws.Range("L137").Value = 107175.669
$ws.Range("M137").Value = -47161.89
$ws.Range("N137").Value = -112275.669

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 35845.766
$ws.Range("I138").Value = 1983.125
$ws.Range("J138").Value = 74545.92999999999
$ws.Range("K138").Value = 5949.375
$ws.Range("L138").Value = 223637.79
$ws.Range("M138").Value = -809.375
$ws.Range("N138").Value = -233917.79

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2822.1667
$ws.Range("I141").Value = 2757
$ws.Range("J141").Value = 2952.5
$ws.Range("K141").Value = 8271
$ws.Range("L141").Value = 8857.5
$ws.Range("M141").Value = -3091
$ws.Range("N141").Value = -19217.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2184.8262
$ws.Range("I2").Value = 2297.0476
$ws.Range("J2").Value = 1006.5
$ws.Range("K2").Value = 2297.0476
$ws.Range("L2").Value = 1006.5
$ws.Range("M2").Value = -2184.0476
$ws.Range("N2").Value = -1232.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20221.25
$ws.Range("I32").Value = 20875.32
$ws.Range("J32").Value = 8666
$ws.Range("K32").Value = 20875.32
$ws.Range("L32").Value = 8666
$ws.Range("M32").Value = -20588.32
$ws.Range("N32").Value = -9240

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 16342
$ws.Range("I33").Value = 16342
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 16342
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -16013

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 44106.332
$ws.Range("I44").Value = 38927.8
$ws.Range("J44").Value = 69999
$ws.Range("K44").Value = 38927.8
$ws.Range("L44").Value = 69999
$ws.Range("M44").Value = -38439.8
$ws.Range("N44").Value = -70975

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2288.04
$ws.Range("I45").Value = 1125.0769
$ws.Range("J45").Value = 3547.9167
$ws.Range("K45").Value = 1125.0769
$ws.Range("L45").Value = 3547.9167
$ws.Range("M45").Value = -748.0769
$ws.Range("N45").Value = -4301.9167

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9373.308000000001
$ws.Range("I61").Value = 1119.125
$ws.Range("J61").Value = 22580
$ws.Range("K61").Value = 1119.125
$ws.Range("L61").Value = 22580
$ws.Range("M61").Value = -907.125
$ws.Range("N61").Value = -23004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 471563.06
$ws.Range("I74").Value = 2000637.4
$ws.Range("J74").Value = 12840.8
$ws.Range("K74").Value = 2000637.4
$ws.Range("L74").Value = 12840.8
$ws.Range("M74").Value = -1999763.4
$ws.Range("N74").Value = -14588.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 471563.06
$ws.Range("I77").Value = 2000637.4
$ws.Range("J77").Value = 12840.8
$ws.Range("K77").Value = 10003187
$ws.Range("L77").Value = 64204
$ws.Range("M77").Value = -9998819
$ws.Range("N77").Value = -72940

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 793.93616
$ws.Range("I97").Value = 677.381
$ws.Range("J97").Value = 1773
$ws.Range("K97").Value = 677.381
$ws.Range("L97").Value = 1773
$ws.Range("M97").Value = -181.381
$ws.Range("N97").Value = -2765

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1500
$ws.Range("I102").Value = 1500
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1500
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 122
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2184.8262
$ws.Range("I116").Value = 2297.0476
$ws.Range("J116").Value = 1006.5
$ws.Range("K116").Value = 2297.0476
$ws.Range("L116").Value = 1006.5
$ws.Range("M116").Value = -3.047599999999875
$ws.Range("N116").Value = -5594.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1360.8889
$ws.Range("I132").Value = 726.9048
$ws.Range("J132").Value = 3579.8333
$ws.Range("K132").Value = 2180.7144
$ws.Range("L132").Value = 10739.4999
$ws.Range("M132").Value = 349.2856000000002
$ws.Range("N132").Value = -15799.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 9373.308000000001
$ws.Range("I136").Value = 1119.125
$ws.Range("J136").Value = 22580
$ws.Range("K136").Value = 3357.375
$ws.Range("L136").Value = 67740
$ws.Range("M136").Value = -807.375
$ws.Range("N136").Value = -72840

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2184.8262
$ws.Range("I3").Value = 2297.0476
$ws.Range("J3").Value = 1006.5
$ws.Range("K3").Value = 2297.0476
$ws.Range("L3").Value = 1006.5
$ws.Range("M3").Value = -2183.0476
$ws.Range("N3").Value = -1234.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 748.64703
$ws.Range("I80").Value = 606.2222
$ws.Range("J80").Value = 908.875
$ws.Range("K80").Value = 606.2222
$ws.Range("L80").Value = 908.875
$ws.Range("M80").Value = 391.7778
$ws.Range("N80").Value = -2904.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 748.64703
$ws.Range("I83").Value = 606.2222
$ws.Range("J83").Value = 908.875
$ws.Range("K83").Value = 3031.111
$ws.Range("L83").Value = 4544.375
$ws.Range("M83").Value = 1960.889
$ws.Range("N83").Value = -14528.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1046.68
$ws.Range("I134").Value = 881.9583
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 2645.8749
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -110.8748999999998
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1568.6666
$ws.Range("I94").Value = 1246.6364
$ws.Range("J94").Value = 1922.9
$ws.Range("K94").Value = 1246.6364
$ws.Range("L94").Value = 1922.9
$ws.Range("M94").Value = -795.6364000000001
$ws.Range("N94").Value = -2824.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6669.727
$ws.Range("I99").Value = 5237.2
$ws.Range("J99").Value = 20995
$ws.Range("K99").Value = 5237.2
$ws.Range("L99").Value = 20995
$ws.Range("M99").Value = -3739.2
$ws.Range("N99").Value = -23991

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 6669.727
$ws.Range("I126").Value = 5237.2
$ws.Range("J126").Value = 20995
$ws.Range("K126").Value = 15711.6
$ws.Range("L126").Value = 62985
$ws.Range("M126").Value = -13241.6
$ws.Range("N126").Value = -67925

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1166.909
$ws.Range("I117").Value = 750
$ws.Range("J117").Value = 1259.5555
$ws.Range("K117").Value = 2250
$ws.Range("L117").Value = 3778.6665
$ws.Range("M117").Value = 1192
$ws.Range("N117").Value = -10662.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 12412.333
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 12412.333
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 12412.333
$ws.Range("N12").Value = -12692.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 39029.5
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 39029.5
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 39029.5
$ws.Range("N49").Value = -39397.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6647.8335
$ws.Range("I70").Value = 7962.6665
$ws.Range("J70").Value = 5333
$ws.Range("K70").Value = 7962.6665
$ws.Range("L70").Value = 5333
$ws.Range("M70").Value = -7692.6665
$ws.Range("N70").Value = -5873

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6647.8335
$ws.Range("I73").Value = 7962.6665
$ws.Range("J73").Value = 5333
$ws.Range("K73").Value = 7962.6665
$ws.Range("L73").Value = 5333
$ws.Range("M73").Value = -7026.6665
$ws.Range("N73").Value = -7205

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 946.8108
$ws.Range("I97").Value = 879.53845
$ws.Range("J97").Value = 1105.8182
$ws.Range("K97").Value = 879.53845
$ws.Range("L97").Value = 1105.8182
$ws.Range("M97").Value = -383.53845
$ws.Range("N97").Value = -2097.8182

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 19184.465
$ws.Range("I102").Value = 21871.459
$ws.Range("J102").Value = 3062.5
$ws.Range("K102").Value = 21871.459
$ws.Range("L102").Value = 3062.5
$ws.Range("M102").Value = -20249.459
$ws.Range("N102").Value = -6306.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1971.2
$ws.Range("I113").Value = 2224.7144
$ws.Range("J113").Value = 1648.5454
$ws.Range("K113").Value = 2224.7144
$ws.Range("L113").Value = 1648.5454
$ws.Range("M113").Value = -54.71439999999984
$ws.Range("N113").Value = -5988.5454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3763.7693
$ws.Range("I126").Value = 1916.6666
$ws.Range("J126").Value = 5347
$ws.Range("K126").Value = 5749.9998
$ws.Range("L126").Value = 16041
$ws.Range("M126").Value = -3279.9998
$ws.Range("N126").Value = -20981

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1680.0476
$ws.Range("I132").Value = 1514.05
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 4542.15
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -2012.15
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2495.6
$ws.Range("I100").Value = 2394.5293
$ws.Range("J100").Value = 3068.3333
$ws.Range("K100").Value = 2394.5293
$ws.Range("L100").Value = 3068.3333
$ws.Range("M100").Value = -1853.5293
$ws.Range("N100").Value = -4150.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2919.2693
$ws.Range("I136").Value = 2812.4736
$ws.Range("J136").Value = 3209.1428
$ws.Range("K136").Value = 8437.4208
$ws.Range("L136").Value = 9627.428400000001
$ws.Range("M136").Value = -5887.4208
$ws.Range("N136").Value = -14727.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 17618.818
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 17618.818
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 17618.818
$ws.Range("N45").Value = -18600.818

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1229.9286
$ws.Range("I113").Value = 1178.6666
$ws.Range("J113").Value = 1322.2
$ws.Range("K113").Value = 3535.9998
$ws.Range("L113").Value = 3966.6
$ws.Range("M113").Value = -1365.9998
$ws.Range("N113").Value = -8306.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 33549.41
$ws.Range("I132").Value = 40537.61
$ws.Range("J132").Value = 2102.5
$ws.Range("K132").Value = 121612.83
$ws.Range("L132").Value = 6307.5
$ws.Range("M132").Value = -119082.83
$ws.Range("N132").Value = -11367.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 26027.48
$ws.Range("I136").Value = 33476.21
$ws.Range("J136").Value = 2439.8333
$ws.Range("K136").Value = 100428.63
$ws.Range("L136").Value = 7319.499899999999
$ws.Range("M136").Value = -97878.63
$ws.Range("N136").Value = -12419.4999
